# Update Name of Algo
# Apply updated RandomForest imputation results to column B (std) and
# column C (mean) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.2973
$ws.Range("B12").Value = 5.026099999999997
$ws.Range("C14").Value = -13.55769999999999
$ws.Range("C26").Value = -11.70849999999999
$ws.Range("B27").Value = 6.493500000000004
$ws.Range("C31").Value = -13.2815
$ws.Range("B32").Value = 6.364499999999998
$ws.Range("C35").Value = -12.43270000000001
$ws.Range("B36").Value = 9.166900000000002
$ws.Range("C37").Value = -13.3138
$ws.Range("B38").Value = 5.014099999999997
$ws.Range("C45").Value = -13.82249999999999
$ws.Range("B46").Value = 7.893500000000006
$ws.Range("C52").Value = -11.0731
$ws.Range("B54").Value = 4.696600000000001
$ws.Range("B55").Value = 5.314199999999997
$ws.Range("B56").Value = 4.651500000000001
$ws.Range("C57").Value = -14.51159999999999
$ws.Range("B67").Value = 5.531099999999992
$ws.Range("B69").Value = 5.560799999999993
$ws.Range("B72").Value = 5.852299999999998
$ws.Range("C81").Value = -13.193
$ws.Range("B83").Value = 5.033799999999996
$ws.Range("C83").Value = -14.60769999999999
$ws.Range("B86").Value = 4.828300000000004
$ws.Range("B91").Value = 5.2111
$ws.Range("B93").Value = 6.838000000000005
$ws.Range("B99").Value = 4.656799999999997
$ws.Range("C100").Value = -12.8156
$ws.Range("C102").Value = -12.7968
